$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header week-ending dates (row 1) ---
# H1 updated, I1 and J1 newly filled in (style already in place, untouched).
$ws.Range("H1").Value = 45632
$ws.Range("I1").Value = 45639
$ws.Range("J1").Value = 45646

# --- New / changed "Person" assignments (column B) ---
# Order matters here: new shared-string entries are interned in first-use
# order, so touch B13/B14/B11 (which introduce brand new strings) before
# B12 (which only reuses an already-existing string).
$ws.Range("B13").Value = "Aidien "
$ws.Range("B14").Value = "Aidien sean "
$ws.Range("B11").Value = "Aidien and Gabe"
$ws.Range("B12").Value = "Aidien and Jabari"

# --- Color-coding / progress cell fills (columns F-J) ---
# Helper: apply a (fill-foreground, fill-background, bold, size) combo to a
# list of cell addresses on the active sheet.
# NOTE: named parameters don't bind reliably in this PowerShell host, so
# the helper is invoked with positional arguments everywhere below.
function Set-ProgressStyle {
    param(
        [string[]] $Cells,
        [double] $ForeColor,
        [double] $BackColor,
        [bool] $Bold,
        [double] $Size
    )
    foreach ($addr in $Cells) {
        $cell = $ws.Range($addr)
        $cell.Interior.Color = $ForeColor
        $cell.Interior.PatternColor = $BackColor
        $cell.Font.Bold = $Bold
        $cell.Font.Size = $Size
        $cell.Font.Name = "Calibri"
    }
}

# style index 26 in the target workbook: fgColor 00B050 / bgColor 00FF00
Set-ProgressStyle @("I3","J3","I4","J4","J5","J6","I7","J7","I9","J9") 5287936 65280 $false 11

# style index 25: fgColor 00B050 / bgColor FFFF00
Set-ProgressStyle @("I5","I6") 5287936 65535 $false 11

# style index 21: fgColor FFFF00 / bgColor FF0000
Set-ProgressStyle @("I8","I10","F27","G27","H27") 65535 255 $false 11

# style index 24: fgColor FFFF00 / bgColor 00FF00
Set-ProgressStyle @("J8","J10","J11","J12","I13","J13","J14") 65535 65280 $false 11

# style index 10: fgColor FFFF00 / bgColor FFFF00
Set-ProgressStyle @("I11","I12","I14","I27") 65535 65535 $false 11

# style index 27 (newly introduced bold, size-9 variant): fgColor FFFF00 / bgColor 00FF00
Set-ProgressStyle @("J27") 65535 65280 $true 9

# --- Restore the active selection to match the saved view ---
[void]$ws.Range("M27").Select()
